$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 7 (HY2 is at row 6; HY3 goes right after it,
# keeping the BrowseProduct column alphabetically sorted), which
# shifts the old rows 7-14 down to 8-15.
$ws.Rows.Item(7).Insert()

# Populate the newly inserted row with the new band-area parameters.
$ws.Cells.Item(7, 1).Value2 = "HY3"
$ws.Cells.Item(7, 2).Value2 = "BA1200"
$ws.Cells.Item(7, 3).Value2 = "BA1450"
$ws.Cells.Item(7, 4).Value2 = "BA1900"

# Update the hidden _FilterDatabase defined name so its range covers
# the new last row (D15 instead of D14).
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$D`$15"
    }
}

# Move the active selection to A8 (matches the saved selection state).
$ws.Range("A8").Select()
